# Word COM-interop script: apply the "human evolution" -> "science" content
# rewrite, fix the misspelled font name, and append the trailing blank
# paragraph, matching the target OOXML diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Fix the misspelled font "TimesNewToman" -> "Times New Roman"
#    everywhere in the document in a single formatting-only
#    find/replace (no text is touched).
# ---------------------------------------------------------------------
$fontFind = $d.Content.Find
$fontFind.ClearFormatting()
$fontFind.Font.Name = "TimesNewToman"
$fontFind.Replacement.ClearFormatting()
$fontFind.Replacement.Font.Name = "Times New Roman"
$fontFind.Execute("", $false, $false, $false, $false, $false, $true, 1, $true, "", 2) | Out-Null

# ---------------------------------------------------------------------
# helper for plain text replacements
# ---------------------------------------------------------------------
function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------
# 2) Title
# ---------------------------------------------------------------------
Replace-Text "Echoes of the Past: Unveiling the Secrets of Human Evolution" "The Beauty and Power of Science: Shaping Our World"

# ---------------------------------------------------------------------
# 3) Byline
# ---------------------------------------------------------------------
Replace-Text " Evelyn Sterling" " Emily Carter"

# ---------------------------------------------------------------------
# 4) Email address line:
#    "evelyn" + "." + "sterling@cambridge-uni" + "." + "edu"
#      -> "emilycarter@galaxyinstitute" + "." + "org"
# ---------------------------------------------------------------------
Replace-Text "evelyn" "emilycarter@galaxyinstitute"
Replace-Text "sterling@cambridge-uni.edu" "org"

# ---------------------------------------------------------------------
# 5) Main body paragraph
# ---------------------------------------------------------------------
Replace-Text "As humans, we are captivated by our origins, yearning to understand the intricate web that weaves us into the tapestry of life" "In the ever-evolving tapestry of human knowledge, science stands as a beacon of progress, illuminating the mysteries of our physical world and empowering us to shape our destiny"

Replace-Text " In this enthralling quest, we delve into the realm of evolutionary science, seeking the echoes of the past that hold the secrets of our existence" " From the intricate mechanics of the cosmos to the profound intricacies of life, the study of science ignites our curiosity, expands our understanding, and opens up a world of possibilities"

Replace-Text " With each discovery, a new chapter of our story unfolds, revealing the remarkable journey that has shaped our species" " This exploration is a journey of discovery, where each step reveals a deeper appreciation for the interconnectedness of all things and inspires us to strive for a better future"

Replace-Text "Unraveling the genetic code, we glimpse into the blueprints of our ancestors, deciphering the intricate language that shapes our traits and characteristics" "Science is woven into the fabric of our daily lives, informing everything from the materials we use to the technologies that connect us"

Replace-Text " Through the analysis of fossilized remains, we piece together the physical evidence of our evolutionary lineage, reconstructing the lives of our predecessors" " It is the foundation upon which countless advancements have been made in medicine, communication, transportation, and countless other fields, improving our quality of life and extending our reach"

Replace-Text " And as we explore the vast expanse of our planet's ecosystems, we uncover the intricate web of interactions that have shaped our species, revealing the interconnectedness of life" " As we continue to unravel the mysteries of the natural world, we unlock the potential for even greater innovation, addressing global challenges, and crafting a more sustainable and harmonious existence"

Replace-Text "The story of human evolution is a tale of adaptation, resilience, and survival" "Moreover, the pursuit of science is a testament to the indomitable spirit of human curiosity and creativity"

Replace-Text " It is a chronicle of countless challenges met and overcome, of victories and setbacks, of triumphs and tragedies" " It is a realm where we can transcend limitations, challenge conventional wisdom, and create new knowledge"

# This also absorbs (and removes) the two trailing runs that followed it
# (the "." and the final "As we continue to explore..." sentence), since
# the replacement text spans all four original runs.
Replace-Text " It is a narrative of the enduring spirit that has propelled our species to the forefront of the animal kingdom, making us the architects of our own destiny. As we continue to explore the depths of our evolutionary history, we unlock the secrets of our past, gaining a profound appreciation for the intricate tapestry of life that has brought us to where we are today." " Each discovery, no matter how small, contributes to a larger tapestry of understanding, empowering us to unravel the enigmas of the universe and to shape a future that is shaped by knowledge, innovation, and boundless curiosity."

# ---------------------------------------------------------------------
# 6) Summary paragraph
# ---------------------------------------------------------------------
Replace-Text "In this exploration of human evolution, we delve into the intricate mechanisms that have shaped our species" "In the vast expanse of human endeavors, science shines as a beacon of progress, offering us the tools to comprehend the world and the ability to shape our destiny"

Replace-Text " Through the analysis of genetic codes, the examination of fossilized remains, and the study of ecological interactions, we unravel the secrets of our past" " It is a pursuit that ignites curiosity, expands understanding, and empowers us to tackle global challenges"

Replace-Text " Our evolutionary journey is a narrative of adaptation, resilience, and survival, revealing the remarkable odyssey that has brought us to the forefront of the animal kingdom" " The study of science is a journey of discovery, innovation, and boundless potential, inspiring us to create a better future"

Replace-Text " As we continue to probe the depths of our evolutionary history, we unlock the mysteries of our origins, gaining a profound appreciation for the interconnectedness of life and the enduring spirit of the human species" " As we continue to unravel the mysteries of the natural world, we unlock the potential for even greater advancements and a world of endless possibilities"

# ---------------------------------------------------------------------
# 7) Append a trailing blank paragraph at the very end of the document
#    (right before the final section break), matching the added <w:p/>.
# ---------------------------------------------------------------------
$d.Content.InsertParagraphAfter()

Write-Host "edit complete"
